$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.824.63"
$ws.Range("E2").Value = "  -0.62%  "

# Row 3
$ws.Range("D3").Value = "3.798.25"
$ws.Range("E3").Value = "  -2.54%  "

# Row 4
$ws.Range("E4").Value = "  -0.25%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "598.61"
$ws.Range("E5").Value = "  -0.34%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.42"
$ws.Range("E6").Value = "  -2.03%  "

# Row 7
$ws.Range("D7").Value = "3.797.78"
$ws.Range("E7").Value = "  -2.67%  "

# Row 8
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("E9").Value = "  -0.49%  "

# Row 10
$ws.Range("E10").Value = "  +0.35%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("E11").Value = "  +1.10%  "

# Row 12
$ws.Range("E12").Value = "  +0.22%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000271"
$ws.Range("E13").Value = "  +4.97%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.82"
$ws.Range("E14").Value = "  -1.00%  "

# Row 15
$ws.Range("D15").Value = "4.441.45"
$ws.Range("E15").Value = "  -2.64%  "

# Row 16
$ws.Range("D16").Value = "3.817.79"
$ws.Range("E16").Value = "  -2.35%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.85"
$ws.Range("E17").Value = "  +3.83%  "

# Row 18
$ws.Range("D18").Value = "67.885.98"
$ws.Range("E18").Value = "  -0.85%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.31"
$ws.Range("E19").Value = "  -1.03%  "

# Row 20
$ws.Range("E20").Value = "  +0.60%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.63"
$ws.Range("E21").Value = "  -2.36%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.83"
$ws.Range("E22").Value = "  -0.10%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.730"
$ws.Range("E23").Value = "  -1.53%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000150"
$ws.Range("E24").Value = "  -6.56%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.54"
$ws.Range("E25").Value = "  +0.10%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.30"
$ws.Range("E26").Value = "  +2.44%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.19"
$ws.Range("E27").Value = "  +0.58%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.31"

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("E29").Value = "  -0.19%  "

# Row 30
$ws.Range("E30").Value = "  -1.17%  "

# Row 31
$ws.Range("D31").Value = "3.960.19"
$ws.Range("E31").Value = "  -2.38%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.69"
$ws.Range("E32").Value = "  -0.73%  "

# Row 33
$ws.Range("E33").Value = "  -2.12%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.53"
$ws.Range("E34").Value = "  -2.70%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.23"
$ws.Range("E35").Value = "  -2.80%  "

# Row 36
$ws.Range("D36").Value = "3.766.46"
$ws.Range("E36").Value = "  -2.80%  "

# Row 37
$ws.Range("B37").Value = "Hedera"
$ws.Range("C37").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.105"
$ws.Range("E37").Value = "  +0.66%  "

# Row 38
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.77"
$ws.Range("E38").Value = "  +1.10%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.92"
$ws.Range("E39").Value = "  +0.00%  "

# Row 40
$ws.Range("B40").Value = "Mantle"
$ws.Range("C40").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.01"
$ws.Range("E40").Value = "  -1.74%  "

# Row 41
$ws.Range("B41").Value = "Kaspa"
$ws.Range("C41").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.138"
$ws.Range("E41").Value = "  -1.73%  "

# Row 42
$ws.Range("E42").Value = "  -0.20%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.318"
$ws.Range("E43").Value = "  +1.15%  "

# Row 45
$ws.Range("E45").Value = "  -0.92%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.74"
$ws.Range("E46").Value = "  +1.09%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "409.76"
$ws.Range("E47").Value = "  -3.79%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "46.33"
$ws.Range("E48").Value = "  -1.99%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000278"
$ws.Range("E49").Value = "  -9.20%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.21"
$ws.Range("E50").Value = "  -1.11%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0355"
$ws.Range("E51").Value = "  -0.93%  "
